$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 259, shifting existing rows 259-273 down to 260-274.
$ws.Rows.Item(259).Insert()

# Populate the newly inserted row 259 with the new data record.
$ws.Cells.Item(259, 1).Value = 4
$ws.Cells.Item(259, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(259, 3).Value = 'Los Lagos'
$ws.Cells.Item(259, 4).Value = 44516
$ws.Cells.Item(259, 4).Style = $ws.Cells.Item(260, 4).Style
$ws.Cells.Item(259, 4).NumberFormat = $ws.Cells.Item(260, 4).NumberFormat
$ws.Cells.Item(259, 5).Value = 10
$ws.Cells.Item(259, 6).Value = 100114001
$ws.Cells.Item(259, 7).Value = 'Papa'
$ws.Cells.Item(259, 8).Value = 'Pehuenche'
$ws.Cells.Item(259, 9).Value = '1a nueva(o)'
$ws.Cells.Item(259, 10).Value = 600
$ws.Cells.Item(259, 11).Value = 15000
$ws.Cells.Item(259, 12).Value = 16000
$ws.Cells.Item(259, 13).Value = 15500
$ws.Cells.Item(259, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(259, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(259, 16).Value = 620
$ws.Cells.Item(259, 17).Value = 25
$ws.Cells.Item(259, 18).Value = 'Hortaliza'
